$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the trailing blank row (old row 12) ---
$ws.Rows.Item(12).Delete()

# --- Remove old row 11 content (ID 10 "Search filter by teacher's name" requirement) ---
$ws.Range("A11:B11").ClearContents()

# --- Header row ---
$ws.Range("A1").Value = "ID_REQUIREMENT"
$ws.Range("C1").Value = "REQUIREMENT_DESCRIPTION"
$ws.Range("B1").Value = "REQUIREMENT_NAME"

# --- Column A (requirement IDs) ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 7
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9

# --- Column B (requirement names), matching final row order ---
$ws.Range("B2").Value = "Implementation of a workflow for ""User Profile"" functionalities"
$ws.Range("B3").Value = "Implementation of a built-in chat system"
$ws.Range("B4").Value = "Integration with a geolocalization provider"
$ws.Range("B5").Value = "Usage of real personal data"
$ws.Range("B6").Value = "Integration with OS Calendar"
$ws.Range("B7").Value = "Implementation of a complete rating system"
$ws.Range("B8").Value = "Implementation of a workflow allowing to save ""connections"" with other people"
$ws.Range("B9").Value = "Implementation of a complex Search Functionality"
$ws.Range("B10").Value = "Implementation of a complex form to publish lessons timetable"

# --- Column C (requirement descriptions) ---
$ws.Range("C2").Value = "The users should be able to compile their own ""User Profile"" form, providing a chosen set of information." + [char]10 + "Every  user should also be able to consult other users' ""User Profile""."
$ws.Range("C9").Value = "A user searching for a private lesson should be able to, at least:" + [char]10 + "1. Search for a teacher directly by name or email;" + [char]10 + "2. Filter the research by field of interest;" + [char]10 + "3. Sort by price, feedback of other users, map location."
$ws.Range("C3").Value = "The users should be able to start an instant messaging communication with other users." + [char]10 + "The users should also be able to access their chats in an easy and structured manner."
$ws.Range("C4").Value = "The teachers should be able to set an indicative geographical location for their private lessons with the help of a geolocalization provider." + [char]10 + "The students should be able to use such geolocalization provider to access the details of the geographical location set by the teacher."
$ws.Range("C5").Value = "The users should be encouraged to use their personal data by connecting their account with a Social Network Active Directory. Facebook is the preferred choice for this."
$ws.Range("C6").Value = "The Mobile App should use OS APIs to connect with the user's default Calendar App to notify them about upcoming private lessons."
$ws.Range("C7").Value = "Both teachers and students should be able to give a public feedback about the persons they are interacting with. The users should be given the possibility to send a feedback after a certain time after a scheduled private lessons and they should be able to, at least, assing an overall rating and a description."
$ws.Range("C8").Value = "The users should be able to flag other users as ""favorite"": this could allow the users to have their favorites marked in every UI of the Mobile App."
$ws.Range("C10").Value = "For every published private lesson offer, a teacher should be able to define and publish a complete timetable through the help of a form."

Write-Host "data done"

# --- Apply the thin border to the whole table (header + body) in one shot ---
$table = $ws.Range("A1:C10")
$tb = $table.Borders
$tb.LineStyle = 1
$tb.Weight = 2

# --- Column A (id) extra style: number format 000 ---
$colA = $ws.Range("A2:A10")
$colA.NumberFormat = "000"

# --- Body: vertical-center alignment ---
$body = $ws.Range("A2:C10")
$body.VerticalAlignment = -4108

# --- Column C (description) extra style: wrap text ---
$colC = $ws.Range("C2:C10")
$colC.WrapText = $true

# --- Header row style: bold white font, blue fill ---
$header = $ws.Range("A1:C1")
$hf = $header.Font
$hf.Bold = $true
$hf.Color = 16777215
$header.Interior.Color = 12611584

Write-Host "style done"

# --- Row heights (wrapped multi-line descriptions) ---
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 30

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 16.1666666667
$ws.Columns.Item(2).ColumnWidth = 73.3072916667
$ws.Columns.Item(3).ColumnWidth = 63.8333333333

# --- Selection ---
$ws.Range("G8").Select()

Write-Host "layout done"
